$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, copying the formatting used by the other headers (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("H1").Value = "Save"

# Populate the new "Save" column data for rows 2-11
$saveValues = @(1, 0, 1, 0, 0, 1, 0, 1, 0, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
